# Add an "Amount" column to the "products" sheet, between "Product Name"
# (column D) and "Pieces in stock" (old column E, now F), plus a matching
# "Amount: In Euros" note row right under the "Internal Id" note row, and
# make "products" the active/selected sheet again (the saved file had
# drifted to "ring_subcategory", the last tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("products")

# 1) Insert the new "Amount" column at E (shifts old E..T right to F..U,
#    carrying their values/styles/widths with them - this also naturally
#    widens the yellow note boxes on rows 16/18/19/20/25 by one column,
#    matching how Excel grows a formatted block when a column is spliced
#    into its interior).
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = 14.109375

# 2) Insert a new note row at 15 (shifts old 15..25 down to 16..26), to
#    host the "Amount: In Euros" explanatory note directly below the
#    "Internal Id" note (row 14).
$ws.Rows("15:15").Insert()

# 3) Populate the new column's header/sample cells.
$ws.Range("E1").Value = "Amount"
$ws.Range("E2").Value = 100

# 4) Populate the new note row, copying row 14's formatting down first so
#    the fill/border/font match the other note rows exactly.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial()
$ws.Range("A15").Value = "Amount: In Euros"

# 5) Restore "products" as the selected/active sheet + selection cell,
#    undoing the drift to "ring_subcategory" (last sheet) recorded in the
#    saved workbook view state.
$ws.Range("E17").Select()
